$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = $ws.Range("B5").Text + ".20189"
$ws.Range("B6").Value = $ws.Range("B6").Text + ".70115"
$ws.Range("B7").Value = $ws.Range("B7").Text + ".20396"

$ws.Range("B8").Select()
